$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 00:22"

# Row 4 - Estados Unidos: refreshed case numbers
$ws.Range("B4").Value = 426300
$ws.Range("C4").Value = 25965
$ws.Range("D4").Value = 22233
$ws.Range("E4").Value = 389445
$ws.Range("F4").Value = 9234
$ws.Range("G4").Value = 1781
$ws.Range("H4").Value = 14622

# Rows 7 & 8 - Germany/France swapped order, with updated numbers
$ws.Range("A7").Value = "Alemania"
$ws.Range("B7").Value = 113067
$ws.Range("C7").Value = 5404
$ws.Range("D7").Value = 46300
$ws.Range("E7").Value = 64511
$ws.Range("F7").Value = 4895
$ws.Range("G7").Value = 240
$ws.Range("H7").Value = 2256

$ws.Range("A8").Value = "Francia"
$ws.Range("B8").Value = 112950
$ws.Range("C8").Value = 3881
$ws.Range("D8").Value = 21254
$ws.Range("E8").Value = 80827
$ws.Range("F8").Value = 7148
$ws.Range("G8").Value = 541
$ws.Range("H8").Value = 10869

# Row 9 - China: updated new-cases and new-deaths figures
$ws.Range("C9").Value = 0
$ws.Range("G9").Value = 0

# Row 87 - Costa Rica: updated figures
$ws.Range("E87").Value = 471
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 2
